$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.74664122137671
$ws.Range("C2").Value = 13.79151025371607
$ws.Range("D2").Value = 14.28263935936553
$ws.Range("E2").Value = 14.69223566107343
$ws.Range("G2").Value = 3.828486085574022
$ws.Range("J2").Value = 8.716142912590508
$ws.Range("K2").Value = 21.37255781607688
$ws.Range("M2").Value = 23.30973710333599
$ws.Range("N2").Value = 29.37594154696663
$ws.Range("B3").Value = 21.59981506691497
$ws.Range("C3").Value = 13.70015538648103
$ws.Range("D3").Value = 14.27207719903124
$ws.Range("E3").Value = 14.70265554667632
$ws.Range("G3").Value = 3.832649061759952
$ws.Range("J3").Value = 8.734211071785962
$ws.Range("K3").Value = 21.28253499357407
$ws.Range("M3").Value = 23.28971347258222
$ws.Range("N3").Value = 29.29338094280023
$ws.Range("B4").Value = 21.51502794665174
$ws.Range("C4").Value = 13.64757390015866
$ws.Range("D4").Value = 14.2685292917703
$ws.Range("E4").Value = 14.71166743632836
$ws.Range("G4").Value = 3.83533548693234
$ws.Range("J4").Value = 8.74606560018209
$ws.Range("K4").Value = 21.23238111370098
$ws.Range("M4").Value = 23.28225802684471
$ws.Range("N4").Value = 29.24356312684595
$ws.Range("B5").Value = 21.48185373781082
$ws.Range("C5").Value = 13.627045962531
$ws.Range("D5").Value = 14.26782271381408
$ws.Range("E5").Value = 14.71599666050414
$ws.Range("G5").Value = 3.836463134494593
$ws.Range("J5").Value = 8.751088119016867
$ws.Range("K5").Value = 21.21324346676085
$ws.Range("M5").Value = 23.28043806546877
$ws.Range("N5").Value = 29.22348579960491
$ws.Range("B6").Value = 21.47642914288739
$ws.Range("C6").Value = 13.62369209233417
$ws.Range("D6").Value = 14.26775004366101
$ws.Range("E6").Value = 14.71675517565814
$ws.Range("G6").Value = 3.83665237092332
$ws.Range("J6").Value = 8.751933695808699
$ws.Range("K6").Value = 21.21014456195473
$ws.Range("M6").Value = 23.28020946274361
$ws.Range("N6").Value = 29.2201656027822
$ws.Range("B7").Value = 21.51457493596139
$ws.Range("C7").Value = 13.64729339066218
$ws.Range("D7").Value = 14.26851676904264
$ws.Range("E7").Value = 14.71172316320181
$ws.Range("G7").Value = 3.835350561367186
$ws.Range("J7").Value = 8.746132558852484
$ws.Range("K7").Value = 21.23211773552961
$ws.Range("M7").Value = 23.2822285486302
$ws.Range("N7").Value = 29.24329144586557
$ws.Range("B8").Value = 21.69491979151089
$ws.Range("C8").Value = 13.75929290206325
$ws.Range("D8").Value = 14.27838838364356
$ws.Range("E8").Value = 14.69528563486463
$ws.Range("G8").Value = 3.829894507266171
$ws.Range("J8").Value = 8.722215215477236
$ws.Range("K8").Value = 21.34046343134702
$ws.Range("M8").Value = 23.30182961779141
$ws.Range("N8").Value = 29.34729547787844
$ws.Range("B9").Value = 22.08978869723615
$ws.Range("C9").Value = 14.00597478717508
$ws.Range("D9").Value = 14.32100910314746
$ws.Range("E9").Value = 14.68381031131648
$ws.Range("G9").Value = 3.820223364270019
$ws.Range("J9").Value = 8.681328712063692
$ws.Range("K9").Value = 21.59291446150259
$ws.Range("M9").Value = 23.37857533999555
$ws.Range("N9").Value = 29.55801482742617
$ws.Range("B10").Value = 22.40306271137748
$ws.Range("C10").Value = 14.20252688217557
$ws.Range("D10").Value = 14.36642109785818
$ws.Range("E10").Value = 14.68804990797959
$ws.Range("G10").Value = 3.813736318385115
$ws.Range("J10").Value = 8.654929507703711
$ws.Range("K10").Value = 21.80179725749029
$ws.Range("M10").Value = 23.4581532701786
$ws.Range("N10").Value = 29.71672145835232
$ws.Range("B11").Value = 22.55015267523665
$ws.Range("C11").Value = 14.29499249460095
$ws.Range("D11").Value = 14.3901146017205
$ws.Range("E11").Value = 14.69272925207799
$ws.Range("G11").Value = 3.810917651095189
$ws.Range("J11").Value = 8.643704567068944
$ws.Range("K11").Value = 21.90166752278225
$ws.Range("M11").Value = 23.49934126712008
$ws.Range("N11").Value = 29.78971228078056
$ws.Range("B12").Value = 22.60646978741118
$ws.Range("C12").Value = 14.33042103114315
$ws.Range("D12").Value = 14.3995200896276
$ws.Range("E12").Value = 14.6948962700068
$ws.Range("G12").Value = 3.809869183348586
$ws.Range("J12").Value = 8.639566304572154
$ws.Range("K12").Value = 21.94016066704753
$ws.Range("M12").Value = 23.51564966054153
$ws.Range("N12").Value = 29.81746109506449
$ws.Range("B13").Value = 22.59431410262784
$ws.Range("C13").Value = 14.32277286724001
$ws.Range("D13").Value = 14.39747523557487
$ws.Range("E13").Value = 14.69441200415297
$ws.Range("G13").Value = 3.810094151062895
$ws.Range("J13").Value = 8.640452561468212
$ws.Range("K13").Value = 21.93184085925164
$ws.Range("M13").Value = 23.51210581769601
$ws.Range("N13").Value = 29.81148014259247
$ws.Range("B14").Value = 22.55477378607107
$ws.Range("C14").Value = 14.29789907540196
$ws.Range("D14").Value = 14.39087973345189
$ws.Range("E14").Value = 14.69289962015244
$ws.Range("G14").Value = 3.81083101495926
$ws.Range("J14").Value = 8.643361859681017
$ws.Range("K14").Value = 21.90482099488453
$ws.Range("M14").Value = 23.50066874716907
$ws.Range("N14").Value = 29.79199304439848
$ws.Range("B15").Value = 22.53063335126204
$ws.Range("C15").Value = 14.28271628864026
$ws.Range("D15").Value = 14.38689611976084
$ws.Range("E15").Value = 14.69202467016397
$ws.Range("G15").Value = 3.811284822983828
$ws.Range("J15").Value = 8.645158512493136
$ws.Range("K15").Value = 21.88835768328406
$ws.Range("M15").Value = 23.4937556696762
$ws.Range("N15").Value = 29.78007063434005
$ws.Range("B16").Value = 22.3935386978249
$ws.Range("C16").Value = 14.19654333642628
$ws.Range("D16").Value = 14.36493348702684
$ws.Range("E16").Value = 14.68779941579142
$ws.Range("G16").Value = 3.813923176686638
$ws.Range("J16").Value = 8.655678830722106
$ws.Range("K16").Value = 21.7953661474114
$ws.Range("M16").Value = 23.45556142140213
$ws.Range("N16").Value = 29.71196685942978
$ws.Range("B17").Value = 22.31057965188815
$ws.Range("C17").Value = 14.1444433197474
$ws.Range("D17").Value = 14.35223541441382
$ws.Range("E17").Value = 14.68591160026201
$ws.Range("G17").Value = 3.815575521404648
$ws.Range("J17").Value = 8.662333278474485
$ws.Range("K17").Value = 21.7395448696135
$ws.Range("M17").Value = 23.43340403909496
$ws.Range("N17").Value = 29.67038703956696
$ws.Range("B18").Value = 22.26329745306879
$ws.Range("C18").Value = 14.11476570617156
$ws.Range("D18").Value = 14.34521774279386
$ws.Range("E18").Value = 14.6850847341616
$ws.Range("G18").Value = 3.816538368796942
$ws.Range("J18").Value = 8.666234575328303
$ws.Range("K18").Value = 21.70789575171966
$ws.Range("M18").Value = 23.42112955333806
$ws.Range("N18").Value = 29.64654653711861
$ws.Range("B19").Value = 22.24736423044721
$ws.Range("C19").Value = 14.10476775591091
$ws.Range("D19").Value = 14.34289087770133
$ws.Range("E19").Value = 14.68484925901424
$ws.Range("G19").Value = 3.816866516697058
$ws.Range("J19").Value = 8.667568180402572
$ws.Range("K19").Value = 21.69725920495052
$ws.Range("M19").Value = 23.41705448575829
$ws.Range("N19").Value = 29.6384876158473
$ws.Range("B20").Value = 22.31936619801692
$ws.Range("C20").Value = 14.14995973149814
$ws.Range("D20").Value = 14.35355757567831
$ws.Range("E20").Value = 14.6860857634791
$ws.Range("G20").Value = 3.815398337600705
$ws.Range("J20").Value = 8.66161726237736
$ws.Range("K20").Value = 21.74543991971495
$ws.Range("M20").Value = 23.43571414179221
$ws.Range("N20").Value = 29.67480554766237
$ws.Range("B21").Value = 22.56637131575649
$ws.Range("C21").Value = 14.30519408507624
$ws.Range("D21").Value = 14.39280526125733
$ws.Range("E21").Value = 14.69333312745913
$ws.Range("G21").Value = 3.810614068234588
$ws.Range("J21").Value = 8.642504281400553
$ws.Range("K21").Value = 21.91273926739181
$ws.Range("M21").Value = 23.50400883522763
$ws.Range("N21").Value = 29.79771396986774
$ws.Range("B22").Value = 22.73138185151432
$ws.Range("C22").Value = 14.40904852061216
$ws.Range("D22").Value = 14.4209795258743
$ws.Range("E22").Value = 14.70037197826416
$ws.Range("G22").Value = 3.807597378617243
$ws.Range("J22").Value = 8.630667697523998
$ws.Range("K22").Value = 22.0259988193822
$ws.Range("M22").Value = 23.55278666543166
$ws.Range("N22").Value = 29.87867307974549
$ws.Range("B23").Value = 22.64299911338298
$ws.Range("C23").Value = 14.35340839916578
$ws.Range("D23").Value = 14.4057126555674
$ws.Range("E23").Value = 14.69640477356132
$ws.Range("G23").Value = 3.809197409860086
$ws.Range("J23").Value = 8.63692531246418
$ws.Range("K23").Value = 21.9651992409224
$ws.Range("M23").Value = 23.52637602850313
$ws.Range("N23").Value = 29.83540776314019
$ws.Range("B24").Value = 22.31539251391618
$ws.Range("C24").Value = 14.14746490209618
$ws.Range("D24").Value = 14.35295894603603
$ws.Range("E24").Value = 14.68600621919486
$ws.Range("G24").Value = 3.815478402231951
$ws.Range("J24").Value = 8.661940737853829
$ws.Range("K24").Value = 21.74277338514079
$ws.Range("M24").Value = 23.43466829839897
$ws.Range("N24").Value = 29.67280773949237
$ws.Range("B25").Value = 21.97874862976956
$ws.Range("C25").Value = 13.93646097150752
$ws.Range("D25").Value = 14.30699538642945
$ws.Range("E25").Value = 14.68468992849578
$ws.Range("G25").Value = 3.822730467928666
$ws.Range("J25").Value = 8.691748419220641
$ws.Range("K25").Value = 21.52043278925841
$ws.Range("M25").Value = 23.35372591562
$ws.Range("N25").Value = 29.5003084655085
